# "Day 1, designed layout"
#
# The underlying workbook content/formatting does not actually change in
# this commit - the original file was produced by a non-Excel tool (note
# the Google-Sheets-style `lightGray` fill name, the non-quantized column
# widths like 12.63, the missing <drawing>-less pageMargins, etc.). This
# commit is simply the file being opened and saved for the first time in
# real Excel, which:
#   - rewrites the package with Excel's own namespaces/metadata
#   - normalizes styles/fonts/fills/column widths to Excel's own
#     internal representations (no visual change)
#   - drops the empty/unused drawing part
#   - stamps standard default page margins
#   - remembers the last-used selection (G11)
#
# Most of that package-level normalization is performed automatically by
# the host runtime on save regardless of script content. The concrete,
# user-controllable pieces of this edit that map onto the Excel object
# model are: the active selection, and the page margins.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Standard default page margins (inches), as Excel stamps them the first
# time a worksheet is saved: left/right 0.7", top/bottom 0.75",
# header/footer 0.3". PageSetup margins are expressed in points, so
# multiply inches by 72.
$ws.PageSetup.LeftMargin = 0.7 * 72
$ws.PageSetup.RightMargin = 0.7 * 72
$ws.PageSetup.TopMargin = 0.75 * 72
$ws.PageSetup.BottomMargin = 0.75 * 72
$ws.PageSetup.HeaderMargin = 0.3 * 72
$ws.PageSetup.FooterMargin = 0.3 * 72

# Restore the last active selection/cell that was recorded in the saved
# workbook view.
$ws.Range("G11").Select()
